$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily records for 2021-09-02 .. 2021-09-09
# Columns: row, date(serial), nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila ab.
$data = @(
    @(367, 44441, 0, 3, 131.3485113835376),
    @(368, 44442, 1, 4, 175.1313485113835),
    @(369, 44443, 1, 5, 218.9141856392294),
    @(370, 44444, 0, 5, 218.9141856392294),
    @(371, 44445, 1, 4, 175.1313485113835),
    @(372, 44446, 0, 3, 131.3485113835376),
    @(373, 44447, 0, 3, 131.3485113835376),
    @(374, 44448, 1, 4, 175.1313485113835)
)

foreach ($rowData in $data) {
    $r = $rowData[0]
    $ws.Cells.Item($r, 1).Value = $rowData[1]
    $ws.Cells.Item($r, 2).Value = $rowData[2]
    $ws.Cells.Item($r, 3).Value = $rowData[3]
    $ws.Cells.Item($r, 4).Value = $rowData[4]
}

# Copy the date-column formatting (border/bold/centered/date numfmt) from the
# last pre-existing row down onto the newly added date cells, matching the
# style already used throughout column A.
$ws.Range("A366").Copy()
$ws.Range("A367:A374").PasteSpecial(-4122)
